# Update "想去人数" (want-to-go count) values for a few events across sheets,
# matching the output generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 41
$wsExhibit.Range("F9").Value = 295
$wsExhibit.Range("F10").Value = 3126

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F4").Value = 3

# Sheet "全部类型" (All types) - combined listing
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 41
$wsAll.Range("F10").Value = 295
$wsAll.Range("F11").Value = 3126
$wsAll.Range("F14").Value = 3
